# Rename the three logo images embedded as InlineShapes in the document's
# first-page header and in the primary / first-page footers.
#
#   first-page header  (header1.xml) : BTec_Logo-Orange   image1.jpg -> image2.jpg
#   primary footer      (footer2.xml): PearsonLogo.png     image2.png -> image1.png
#   first-page footer   (footer1.xml): PearsonLogo.png     image2.png -> image1.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Primary footer (wdHeaderFooterPrimary) -> footer2.xml, PearsonLogo, id=2 ---
$primaryFooter = $sec.Footers.Item(1)
if ($primaryFooter.Exists -and $primaryFooter.Range.InlineShapes.Count -ge 1) {
    $logo = $primaryFooter.Range.InlineShapes.Item(1)
    if ($logo.AlternativeText -like "*PearsonLogo*") {
        $logo.Name = "image1.png"
    }
}

# --- First-page footer (wdHeaderFooterFirstPage) -> footer1.xml, PearsonLogo, id=3 ---
$firstPageFooter = $sec.Footers.Item(2)
if ($firstPageFooter.Exists -and $firstPageFooter.Range.InlineShapes.Count -ge 1) {
    $logo = $firstPageFooter.Range.InlineShapes.Item(1)
    if ($logo.AlternativeText -like "*PearsonLogo*") {
        $logo.Name = "image1.png"
    }
}

# --- First-page header (wdHeaderFooterFirstPage) -> header1.xml, BTec_Logo-Orange, id=1 ---
$firstPageHeader = $sec.Headers.Item(2)
if ($firstPageHeader.Exists -and $firstPageHeader.Range.InlineShapes.Count -ge 1) {
    $logo = $firstPageHeader.Range.InlineShapes.Item(1)
    if ($logo.AlternativeText -like "*BTec_Logo*") {
        $logo.Name = "image2.jpg"
    }
}

Write-Output "Renamed header/footer logo InlineShapes"
